# Add files via upload
# Sets cell B5 on the active sheet to "個人" + newline + "トライアル"
# (mirrors the existing "チーム\nトライアル" label already used for the team section in B23).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NASA トライアル")

$ws.Range("B5").Value = "個人" + [char]10 + "トライアル"
